$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","D","E","F","G","I","J","L","N","O")

$data = @(
    @(16.72108471310284, 9.35867802179188, 14.24662713491184, 31.91294490784457, 3.637742375804281, 24.220537026416, 10.41762639207613, 11.74661469197533, 18.17077857607983, 24.20472263912309),
    @(16.40463245790677, 9.365936399759541, 14.22967815492202, 31.90783199814059, 3.640268678967709, 24.36803892057443, 10.40962106939383, 11.57250471215913, 18.18600810926315, 24.17600803214458),
    @(16.20902157080354, 9.370846538412332, 14.21983248234154, 31.91342046326498, 3.64190373001763, 24.4631875786989, 10.40505581860722, 11.4655926140586, 18.19736232951449, 24.16482822832075),
    @(16.12907920341724, 9.372961755673117, 14.21596317411936, 31.9178917533174, 3.64259119073404, 24.50311642839057, 10.40328432914306, 11.42207100980546, 18.20249433243099, 24.1618974478293),
    @(16.11579393977786, 9.373319896956994, 14.2153293600262, 31.91876666470032, 3.642706623358643, 24.50981642644315, 10.40299556821508, 11.41484846216124, 18.20337704056011, 24.16150900186939),
    @(16.20794423536382, 9.370874601822031, 14.21977971860766, 31.91347188361306, 3.641912915568799, 24.46372139281137, 10.40503156664351, 11.46500541692024, 18.19742949477831, 24.164782120162),
    @(16.61230234153771, 9.361086754719906, 14.24066776265073, 31.90937053299225, 3.638596075290031, 24.27044642265573, 10.4147940247566, 11.68661182849413, 18.1756147646319, 24.19348508184428),
    @(17.39064131640431, 9.345479172863172, 14.28600221576644, 31.9705461702714, 3.632754249200177, 23.92765883708097, 10.43667814560018, 12.11904027427767, 18.14867150879975, 24.30076884528174),
    @(17.94816603554842, 9.336183066874186, 14.32186997139407, 32.05754860524618, 3.628861734808696, 23.69771011214397, 10.45437707145627, 12.43271510109984, 18.13844941049821, 24.41033648951314),
    @(18.19762331348306, 9.332422106021411, 14.3387199973032, 32.10619640035746, 3.62717673693952, 23.59781595233022, 10.46276984060666, 12.57398423018243, 18.13585968796686, 24.46675476855033),
    @(18.29140489099835, 9.33106493889248, 14.34517531939153, 32.12591407901752, 3.62655092827763, 23.56066315941359, 10.46599608577962, 12.62723020940001, 18.13517365698951, 24.48905338776078),
    @(18.27123898504234, 9.331354252116858, 14.34378176781608, 32.12161004451441, 3.626685162962233, 23.56863469747929, 10.46529913384824, 12.61577453137822, 18.13530832852062, 24.48420962519704),
    @(18.20535292049014, 9.332309109196055, 14.33924959822853, 32.10779266985777, 3.627125005870414, 23.5947458563253, 10.46303430789342, 12.57837006180457, 18.13579735233839, 24.46857063545103),
    @(18.16490452252898, 9.332902708098958, 14.33648316195832, 32.09949760329841, 3.627396017528941, 23.61082752298604, 10.46165326841546, 12.55542493904545, 18.1361352148507, 24.45911260935419),
    @(17.93177236201516, 9.336438247822441, 14.32077933728591, 32.05455118546134, 3.628973572944759, 23.70433301237676, 10.45383536869515, 12.42345031360223, 18.13865997617689, 24.40678080298529),
    @(17.78762579957344, 9.338726842538126, 14.31128070215553, 32.02929557960588, 3.629963263077892, 23.76290037419377, 10.44912607603466, 12.34209114101022, 18.14073539827393, 24.37635346718027),
    @(17.70432932329107, 9.340087236421333, 14.30586769491742, 32.01562354167245, 3.630540579434418, 23.79703030480532, 10.44644961501982, 12.29516412751432, 18.14212315883873, 24.35947214346686),
    @(17.67606270096973, 9.340555416586277, 14.30404365724604, 32.01114138406787, 3.630737437294224, 23.80866236488579, 10.44554897218537, 12.2792542818565, 18.14262640489697, 24.35386316818644),
    @(17.80301114193614, 9.338478660078007, 14.31228665182915, 32.0318957189169, 3.629857073836671, 23.75661989224679, 10.44962406303067, 12.35076589108727, 18.14049439605245, 24.37952845228729),
    @(18.22472442017373, 9.332026827454518, 14.34057879945678, 32.11181607521321, 3.626995481026269, 23.58705807704988, 10.46369824530433, 12.5893637838316, 18.13564573107062, 24.47313892773495),
    @(18.49632506286171, 9.328200749348253, 14.35950305584768, 32.17159655710601, 3.62519671623308, 23.48017269808142, 10.47317633491317, 12.7438270345019, 18.13419343733656, 24.53975804179972),
    @(18.35176030519523, 9.330207146210549, 14.34936386017345, 32.13900317248186, 3.626150234024027, 23.53686035296955, 10.46809243765744, 12.66153636333189, 18.13481204237362, 24.50370856504761),
    @(17.79605674502836, 9.33859072429305, 14.31183171269141, 32.03071755459188, 3.62990505603929, 23.7594578688669, 10.44939882656227, 12.34684450825685, 18.14060274694707, 24.37809113504699),
    @(17.18221826853436, 9.349319100366426, 14.27328059174956, 31.94659421870041, 3.634264150184959, 24.01653327938501, 10.43046963883295, 12.00257167998925, 18.15427318326629, 24.26631577149783)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range("$($cols[$j])$rowNum").Value = $rowValues[$j]
    }
}
